$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the diff (Test Candidate details for P2/P3)
$ws.Range("A2").Value = "fSgDQ134"
$ws.Range("B2").Value = 23081404
$ws.Range("C2").Value = "rkowlyw64"
$ws.Range("D2").Value = "dMvH`$2&6"
$ws.Range("F2").Value = "DAeTNpEK"
$ws.Range("G2").Value = "Rmvc"
